# Switch to using AB instead of PA
# Update the "PA" column header to "AB" on every sheet, tweak a handful of
# underlying stat values that were corrected at the same time, and restore
# each sheet's last-used selection.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Fall 2015 09.16" ---------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Value = "AB"
$ws1.Range("J4").Value = 1
$ws1.Range("B5").Value = 3

# --- Sheet 2: "Fall 2015 09.09" ---------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = "AB"
$ws2.Range("B2").Value = 4
$ws2.Range("B4").Value = 5
$ws2.Range("B5").Value = 5

# --- Sheet 3: "Spring 2014 04.16" -------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B1").Value = "AB"
$ws3.Range("B2").Value = 4
$ws3.Range("B4").Value = 4

# --- Sheet 4: "Spring 2014 04.09" -------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B1").Value = "AB"
$ws4.Range("B2").Value = 4

# --- Sheet 5: "Tournament Fall 2015" ----------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B1").Value = "AB"
$ws5.Range("B2").Formula = "=4+3+3"
$ws5.Range("B3").Formula = "=3+2+3"
$ws5.Range("B5").Formula = "=4+3+3"
$ws5.Range("B10").Formula = "=3+1+3"

# --- Restore per-sheet selections (in tab order) ----------------------
$ws2.Activate()
$ws2.Range("B1:B6").Select()

$ws3.Activate()
$ws3.Range("J5").Select()

$ws4.Activate()
$ws4.Range("B1:B6").Select()

$ws5.Activate()
$ws5.Range("B1:B13").Select()

# Leave sheet 1 as the active/selected tab, matching the saved workbook.
$ws1.Activate()
$ws1.Range("J4").Select()
